$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.172.08"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "2.480.24"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.72%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").Value = "2.480.14"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("E10").Value = "  +3.40%  "

$ws.Range("E11").Value = "  +0.95%  "

$ws.Range("E12").Value = "  +0.81%  "

$ws.Range("E13").Value = "  +0.46%  "

$ws.Range("D14").Value = "2.941.72"
$ws.Range("E14").Value = "  +0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").Value = "67.052.74"
$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("D18").Value = "2.534.61"
$ws.Range("E18").Value = "  +3.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("E23").Value = "  -0.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.02%  "

$ws.Range("D28").Value = "2.604.53"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.59%  "

$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.14"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.53%  "

$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.64%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("D47").Value = "0.0₆0263"
$ws.Range("E47").Value = "  +4.30%  "

$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("E51").Value = "  -0.34%  "
